# "lots of discharge data" -- add a third ("new depth") discharge table to
# sheet "stn3", recomputed from the recalibrated depth data already on the
# sheet, and leave the workbook with stn3 active / F37 selected (matching the
# author's final view after adding the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 35: bold section label -------------------------------------------------
$ws.Cells.Item(35, 1).Value = "new depth"
$ws.Cells.Item(35, 1).Font.Bold = $true

# --- Row 36: column headers (same labels used by the other two tables) ---------
$ws.Cells.Item(36, 1).Value = "X"
$ws.Cells.Item(36, 2).Value = "V"
$ws.Cells.Item(36, 3).Value = "D"
$ws.Cells.Item(36, 4).Value = "segment"
$ws.Cells.Item(36, 5).Value = "Q"
$ws.Cells.Item(36, 6).Value = "Qtotal"

# --- Rows 37-50: data ------------------------------------------------------------
# Column A: station depth (same grid as the other two tables).
# Column B: velocity, carried over (as values) from the B20:B33 table above.
# Column C: depth re-expressed in cm, pulled from the row 17-above C cell.
# Column D: mirrors column A.
# Column E: incremental discharge segment; first row (37) has none.
# F37: grand total across the new table.
$A = @(0.65, 0.7, 0.75, 0.8, 0.85, 0.9, 0.95, 1, 1.05, 1.1000000000000001, 1.1499999999999999, 1.2, 1.25, 1.3)
$B = @(0, 0.051480000000000005, 0.16016, 0.52623999999999993, 0.58343999999999996, 0.67496, 0.68640000000000001, 0.58916000000000002, 0.41755999999999999, 0.38324000000000003, 0.26884000000000002, 0.13727999999999999, 0.1144, 0)

for ($i = 0; $i -lt 14; $i++) {
    $row = 37 + $i
    $srcRow = $row - 17

    $ws.Cells.Item($row, 1).Value = $A[$i]
    $ws.Cells.Item($row, 2).Value = $B[$i]
    $ws.Cells.Item($row, 3).Formula = "=C$srcRow*2.54"
    $ws.Cells.Item($row, 4).Formula = "=A$row"

    if ($row -gt 37) {
        $prevRow = $row - 1
        $ws.Cells.Item($row, 5).Formula = "=(D$row-D$prevRow)*(B$row)*C$row"
    }
}

$ws.Cells.Item(37, 6).Formula = "=SUM(E37:E55)"

# --- Final view state: stn3 active, F37 selected --------------------------------
$ws.Activate() | Out-Null
$ws.Range("F37").Select() | Out-Null
